# Updates cryptos list values (price & 1h volume change) to the latest scrape,
# and swaps the Polkadot/Litecoin row positions (rows 14 and 15).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value. Numeric-looking Price values are prefixed with a literal
# leading apostrophe so Excel stores them as text (quote-prefix) instead of
# silently converting them to floating point numbers and losing exact digits.
$updates = @{
    'D2' = '27.104.68'
    'E2' = '  +0.38%  '
    'D3' = '1.826.14'
    'E3' = '  +0.13%  '
    'E4' = '  +0.31%  '
    'D5' = '''312.19'
    'E5' = '  +0.30%  '
    'D6' = '''1.006'
    'E6' = '  +0.24%  '
    'D7' = '''0.4694'
    'E7' = '  +0.05%  '
    'D8' = '''0.3655'
    'E8' = '  -0.19%  '
    'D9' = '''0.07385'
    'E9' = '  +0.38%  '
    'D10' = '''0.8776'
    'D11' = '''20.26'
    'E11' = '  -0.23%  '
    'D12' = '1.901.01'
    'E12' = '  +2.98%  '
    'D13' = '''0.07609'
    'E13' = '  +4.09%  '
    'B14' = 'Litecoin'
    'C14' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'D14' = '''93.35'
    'E14' = '  +1.66%  '
    'B15' = 'Polkadot'
    'C15' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'D15' = '''5.367'
    'E15' = '  -1.20%  '
    'D16' = '''6.524'
    'E16' = '  +0.02%  '
    'D17' = '''1.005'
    'E17' = '  +0.00%  '
    'D18' = '''0.000008718'
    'E18' = '  -0.34%  '
    'E19' = '  +0.28%  '
    'D20' = '27.582.99'
    'E20' = '  +2.09%  '
    'D21' = '''14.59'
    'E21' = '  -0.71%  '
    'D22' = '''5.236'
    'E22' = '  -1.01%  '
    'D24' = '2.083.24'
    'E24' = '  +0.85%  '
    'D25' = '''1.871'
    'E25' = '  -1.10%  '
    'D26' = '''150.99'
    'E26' = '  -0.02%  '
    'D27' = '''18.37'
    'E27' = '  +0.05%  '
    'D28' = '''2.132'
    'E28' = '  -0.56%  '
    'D29' = '''5.167'
    'E29' = '  -1.71%  '
    'D30' = '''116.31'
    'E30' = '  -0.29%  '
    'E31' = '  +0.27%  '
    'D32' = '''0.7438'
    'E32' = '  -1.55%  '
    'D33' = '''1.161'
    'E33' = '  -0.19%  '
    'D34' = '''4.509'
    'E34' = '  -0.04%  '
    'D35' = '''2.940'
    'E35' = '  +0.42%  '
    'D36' = '''2.648'
    'E36' = '  +11.57%  '
    'E37' = '  +0.22%  '
    'D38' = '''1.089'
    'E38' = '  -0.63%  '
    'D39' = '''0.05291'
    'E39' = '  -0.39%  '
    'D40' = '''0.01932'
    'E40' = '  -0.91%  '
    'D41' = '''7.306'
    'E41' = '  +1.32%  '
    'D42' = '''2.924'
    'E42' = '  -1.87%  '
    'D43' = '''0.5261'
    'E43' = '  -0.89%  '
    'D44' = '''0.1642'
    'E44' = '  -0.77%  '
    'D45' = '''8.369'
    'E45' = '  -1.39%  '
    'D46' = '''0.4900'
    'E46' = '  +0.07%  '
    'E47' = '  -0.95%  '
    'D48' = '''1.006'
    'E48' = '  +0.26%  '
    'D49' = '''104.38'
    'E49' = '  +1.10%  '
    'D50' = '''1.651'
    'E50' = '  -0.73%  '
    'D51' = '''0.06272'
    'E51' = '  -0.44%  '
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
